# Generate Report for handoff
#
# The localization run for 140339dc-a6bb-4d5a-9aad-9e8f245d665d.md failed the
# handoff transform and a *new* handoff was generated under a new id
# (e43df541-33d8-4984-baae-a2a78290d9c6.md). This updates the report:
#   - the "File Name" display text switches to the new id (the underlying
#     link target is left as-is, only what's shown changes)
#   - status flips from "Ready for handoff" to "Handoff transform failed"
#   - since the transform failed there's no handoff .xlf produced this run,
#     so the per-language "Latest Handoff File" / "Latest Handoff Datetime"
#     cells are cleared back to their not-yet-handed-off defaults and the
#     dependency state resets from "Include" to "Ignored"

$wb = $excel.ActiveWorkbook

$oldMdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/8d204191332c208e3696a339bfc5a92546b586d4/e2e/140339dc-a6bb-4d5a-9aad-9e8f245d665d.md"
$configUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/5fa016d8b9736f78b6dbe292699dda89c1503617/.localization-config"
$newMdDisplay = "e43df541-33d8-4984-baae-a2a78290d9c6.md"
$newStatus    = "Handoff transform failed"
$epoch        = "0001-01-01 00:00:00"

# ---- Overview sheet: just the display text + status text change ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), $oldMdUrl, "", "", $newMdDisplay) | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), $configUrl, "", "", ".localization-config") | Out-Null

# ---- Per-language sheets: zh-cn / de-de ----
foreach ($name in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("B2").Value = $newStatus

    # The handoff transform never produced a target file this run, so the
    # handoff-file link/cell is gone and the handoff datetime resets.
    $ws.Range("C2").ClearContents()
    $ws.Range("D2").Value = $epoch
    $ws.Range("G2").Value = $epoch
    $ws.Range("H2").Value = "Ignored"

    $ws.Range("D3").Value = $epoch
    $ws.Range("G3").Value = $epoch
    $ws.Range("H3").Value = "Ignored"

    # Rebuild the two remaining hyperlinks (md display text changes, the
    # .localization-config one is untouched) -- this also drops the C2
    # handoff-file hyperlink that ClearContents alone wouldn't remove.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $oldMdUrl, "", "", $newMdDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, "", "", ".localization-config") | Out-Null
}
